$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Menu"
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 21
$ws.Range("A3").Value = "Combobox"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 11
$ws.Range("A4").Value = "DataGrid"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 11
$ws.Range("A5").Value = "Tree"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 9
$ws.Range("A6").Value = "Nav"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 9
$ws.Range("A7").Value = "Popover"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 8
$ws.Range("A8").Value = "Dialog"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 6
$ws.Range("A9").Value = "Dropdown"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 6
$ws.Range("A10").Value = "Table"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 5
$ws.Range("A11").Value = "Tooltip"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 5
$ws.Range("A12").Value = "TagPicker"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 4
$ws.Range("A13").Value = "Toolbar"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 4
$ws.Range("A14").Value = "Virtualizer"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 4
$ws.Range("A15").Value = "Skeleton"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 3
$ws.Range("A16").Value = "MessageBar"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 3
$ws.Range("A17").Value = "Calendar Compat"
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 3
$ws.Range("A18").Value = "TeachingPopover"
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 3
$ws.Range("A19").Value = "Toast"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 2
$ws.Range("A20").Value = "Drawer"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 2
$ws.Range("A21").Value = "Slider"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 2
$ws.Range("A22").Value = "List"
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 2
$ws.Range("A23").Value = "FluentProvider"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("A24").Value = "Portal"
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 2
$ws.Range("A25").Value = "Tabs"
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 2
$ws.Range("A26").Value = "Accordion"
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 2
$ws.Range("A27").Value = "Switch"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 2
$ws.Range("A28").Value = "DatePicker"
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("A29").Value = "Image"
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("A30").Value = "Checkbox"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("A31").Value = "Button"
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 1
$ws.Range("A32").Value = "Spinner"
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 1
$ws.Range("A33").Value = "DatePickerCompat"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("A34").Value = "FocusTrapZone"
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 1
$ws.Range("A35").Value = "SplitButton"
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 1
$ws.Range("A36").Value = "InfoLabel"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 1
$ws.Range("A37").Value = "SearchBox"
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 1
$ws.Range("A38").Value = "Input"
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("A39").Value = "Badge"
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 1
$ws.Range("A40").Value = "Label"
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 1
$ws.Range("A41").Value = "Popup"
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 1
$ws.Range("A42").Value = "Carousel"
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 1
$ws.Range("A43").Value = "MenuItem"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 1
$ws.Range("A44").Value = "AvatarGroup"
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 1
$ws.Range("A45").Value = "Avatar"
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 1
$ws.Range("A46").Value = "Tag"
$ws.Range("B46").Value = 0
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("A47").Value = "Pickers"
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("A48").Value = "Rating"
$ws.Range("B48").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("A49").Value = "Keytip"
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("A50").Value = "Segment"
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("A51").Value = "SpinButton"
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("A52").Value = "ColorPicker"
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("A53").Value = "FloatingLabelInput"
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 0

$ws.Rows.Item(54).Delete()
